$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 44 ("pubmed articles"): Test Done / Name go from "na" to "?"
$ws.Range("C44").Value = "?"
$ws.Range("D44").Value = "?"
$ws.Range("E44").Value = "HOW TO INCLUDE PUBMED USING REQUESTS MODULE"

# Row 45 ("pubmed search term"): update comment
$ws.Range("E45").Value = "HOW TO INCLUDE PUBMED USING REQUESTS MODULE"

# Update the active selection to reflect where the author left off editing
$ws.Range("E38").Select()
